# Update the "Pais" (countries COVID stats) worksheet with the newer snapshot
# of data, re-sorted by total cases descending. Columns are:
#   A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#   E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 19:05"

# Estados Unidos (row 4) - refreshed totals, same rank
$ws.Range("B4").Value = 1602132
$ws.Range("C4").Value = 9409
$ws.Range("D4").Value = 371534
$ws.Range("E4").Value = 1135266
$ws.Range("G4").Value = 396
$ws.Range("H4").Value = 95332

# India (row 14) - refreshed totals, same rank
$ws.Range("B14").Value = 115572
$ws.Range("C14").Value = 3544
$ws.Range("D14").Value = 46873
$ws.Range("E14").Value = 65197
$ws.Range("G14").Value = 68
$ws.Range("H14").Value = 3502

# Irak overtakes Azerbaiyan and Camerun (rows 69-71)
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 3877
$ws.Range("C69").Value = 153
$ws.Range("D69").Value = 2483
$ws.Range("E69").Value = 1254
$ws.Range("G69").Value = 6
$ws.Range("H69").Value = 140

$ws.Range("A70").Value = "Azerbaiyan"
$ws.Range("B70").Value = 3749
$ws.Range("C70").Value = 118
$ws.Range("D70").Value = 2340
$ws.Range("E70").Value = 1365
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 44

$ws.Range("A71").Value = "Camerun"
$ws.Range("B71").Value = 3733
$ws.Range("D71").Value = 1595
$ws.Range("E71").Value = 1992
$ws.Range("H71").Value = 146

# Republica de Yibuti overtakes Cuba, Republica de Macedonia
# and Consejo Danes para los Refugiados (rows 86-89)
$ws.Range("A86").Value = "Republica de Yibuti"
$ws.Range("B86").Value = 2047
$ws.Range("C86").Value = 219
$ws.Range("D86").Value = 1055
$ws.Range("E86").Value = 982
$ws.Range("H86").Value = 10

$ws.Range("A87").Value = "Cuba"
$ws.Range("B87").Value = 1908
$ws.Range("C87").Value = 8
$ws.Range("D87").Value = 1603
$ws.Range("E87").Value = 225
$ws.Range("H87").Value = 80

$ws.Range("A88").Value = "Republica de Macedonia"
$ws.Range("B88").Value = 1898
$ws.Range("C88").Value = 40
$ws.Range("D88").Value = 1378
$ws.Range("E88").Value = 409
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 111

$ws.Range("A89").Value = "Consejo Danes para los Refugiados"
$ws.Range("B89").Value = 1835
$ws.Range("C89").Value = 104
$ws.Range("D89").Value = 303
$ws.Range("E89").Value = 1471
$ws.Range("H89").Value = 61

# Zambia overtakes Paraguay (rows 115-116)
$ws.Range("A115").Value = "Zambia"
$ws.Range("B115").Value = 866
$ws.Range("C115").Value = 34
$ws.Range("D115").Value = 302
$ws.Range("E115").Value = 557
$ws.Range("H115").Value = 7

$ws.Range("A116").Value = "Paraguay"
$ws.Range("B116").Value = 836
$ws.Range("C116").Value = 3
$ws.Range("D116").Value = 256
$ws.Range("E116").Value = 569
$ws.Range("H116").Value = 11

# Jordania (row 123) - refreshed totals, same rank
$ws.Range("B123").Value = 684
$ws.Range("C123").Value = 12
$ws.Range("D123").Value = 457
$ws.Range("E123").Value = 218

# Cabo Verde (row 140) - refreshed totals, same rank
$ws.Range("B140").Value = 356
$ws.Range("C140").Value = 7
$ws.Range("E140").Value = 268
